$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 13 used to hold the "sum" total; clear it first so the new
#     data row written into r13 below doesn't inherit/collide with it ---
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()

# --- New data rows 6-8: text placeholders in the date column (keeps the
#     column's default centered style, no date number format) ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 10000
$ws.Range("C6").Value = "Thang Duong"
$ws.Range("D6").Value = "….."
$ws.Range("E6").Value = "late for meeting"
$ws.Range("F6").Value = "Waiting"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 10000
$ws.Range("C7").Value = "Thi Nguyen"
$ws.Range("D7").Value = "…"
$ws.Range("E7").Value = "late for meeting"
$ws.Range("F7").Value = "Waiting"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 10000
$ws.Range("C8").Value = "Thang Duong"
$ws.Range("D8").Value = "…"
$ws.Range("E8").Value = "late for meeting"
$ws.Range("F8").Value = "Waiting"

# --- New data rows 9-14: real dates, formatted like the existing date
#     cells (apply the number format before assigning the value so no
#     extra style/numFmt gets minted) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 10000
$ws.Range("C9").Value = "Kim Hoang"
$ws.Range("D9").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("D9").Value = "9/12/2010"
$ws.Range("E9").Value = "no submit"
$ws.Range("F9").Value = "Waiting"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 10000
$ws.Range("C10").Value = "Thi Nguyen"
$ws.Range("D10").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("D10").Value = "9/12/2010"
$ws.Range("E10").Value = "no submit"
$ws.Range("F10").Value = "Waiting"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 10000
$ws.Range("C11").Value = "Kim Hoang"
$ws.Range("D11").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("D11").Value = "12/20/2010"
$ws.Range("E11").Value = "late for meeting"
$ws.Range("F11").Value = "Waiting"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 10000
$ws.Range("C12").Value = "Thi Nguyen"
$ws.Range("D12").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("D12").Value = "12/20/2010"
$ws.Range("E12").Value = "rule 19"
$ws.Range("F12").Value = "Waiting"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 10000
$ws.Range("C13").Value = "Tan Nguyen"
$ws.Range("D13").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("D13").Value = "12/20/2010"
$ws.Range("E13").Value = "rule 19"
$ws.Range("F13").Value = "Waiting"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 10000
$ws.Range("C14").Value = "Thang Le"
$ws.Range("D14").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("D14").Value = "12/20/2010"
$ws.Range("E14").Value = "rule 19"
$ws.Range("F14").Value = "Waiting"

# --- Re-create the "sum" row at 20, extending the total range ---
$ws.Range("B20").Value = "sum"
$ws.Range("C20").Formula = "=SUM(B2:B16)"

# --- Final selection matches the author's last-saved cursor position ---
$ws.Range("G21").Select() | Out-Null
